$d = $word.ActiveDocument

$found = $d.Content.Find.Execute("21÷8=2, 5", $true, $false, $false, $false, $false, $true, 1, $false, "52÷4=13, 0", 2)
if (-not $found) { throw "Replacement failed: '21÷8=2, 5' -> '52÷4=13, 0'" }
$found = $d.Content.Find.Execute("20÷5=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "58÷7=8, 2", 2)
if (-not $found) { throw "Replacement failed: '20÷5=4, 0' -> '58÷7=8, 2'" }
$found = $d.Content.Find.Execute("69÷4=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "21÷9=2, 3", 2)
if (-not $found) { throw "Replacement failed: '69÷4=17, 1' -> '21÷9=2, 3'" }
$found = $d.Content.Find.Execute("81÷8=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "33÷4=8, 1", 2)
if (-not $found) { throw "Replacement failed: '81÷8=10, 1' -> '33÷4=8, 1'" }
$found = $d.Content.Find.Execute("74÷7=10, 4", $true, $false, $false, $false, $false, $true, 1, $false, "20÷3=6, 2", 2)
if (-not $found) { throw "Replacement failed: '74÷7=10, 4' -> '20÷3=6, 2'" }
$found = $d.Content.Find.Execute("16÷4=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "46÷5=9, 1", 2)
if (-not $found) { throw "Replacement failed: '16÷4=4, 0' -> '46÷5=9, 1'" }
$found = $d.Content.Find.Execute("87÷4=21, 3", $true, $false, $false, $false, $false, $true, 1, $false, "86÷4=21, 2", 2)
if (-not $found) { throw "Replacement failed: '87÷4=21, 3' -> '86÷4=21, 2'" }
$found = $d.Content.Find.Execute("15÷3=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "22÷9=2, 4", 2)
if (-not $found) { throw "Replacement failed: '15÷3=5, 0' -> '22÷9=2, 4'" }
$found = $d.Content.Find.Execute("56÷5=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "64÷5=12, 4", 2)
if (-not $found) { throw "Replacement failed: '56÷5=11, 1' -> '64÷5=12, 4'" }
$found = $d.Content.Find.Execute("69÷7=9, 6", $true, $false, $false, $false, $false, $true, 1, $false, "76÷7=10, 6", 2)
if (-not $found) { throw "Replacement failed: '69÷7=9, 6' -> '76÷7=10, 6'" }
$found = $d.Content.Find.Execute("12÷3=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "79÷2=39, 1", 2)
if (-not $found) { throw "Replacement failed: '12÷3=4, 0' -> '79÷2=39, 1'" }
$found = $d.Content.Find.Execute("67÷4=16, 3", $true, $false, $false, $false, $false, $true, 1, $false, "55÷3=18, 1", 2)
if (-not $found) { throw "Replacement failed: '67÷4=16, 3' -> '55÷3=18, 1'" }
$found = $d.Content.Find.Execute("53÷8=6, 5", $true, $false, $false, $false, $false, $true, 1, $false, "61÷7=8, 5", 2)
if (-not $found) { throw "Replacement failed: '53÷8=6, 5' -> '61÷7=8, 5'" }
$found = $d.Content.Find.Execute("66÷3=22, 0", $true, $false, $false, $false, $false, $true, 1, $false, "63÷3=21, 0", 2)
if (-not $found) { throw "Replacement failed: '66÷3=22, 0' -> '63÷3=21, 0'" }
$found = $d.Content.Find.Execute("33÷7=4, 5", $true, $false, $false, $false, $false, $true, 1, $false, "45÷8=5, 5", 2)
if (-not $found) { throw "Replacement failed: '33÷7=4, 5' -> '45÷8=5, 5'" }
$found = $d.Content.Find.Execute("25÷5=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "40÷8=5, 0", 2)
if (-not $found) { throw "Replacement failed: '25÷5=5, 0' -> '40÷8=5, 0'" }
$found = $d.Content.Find.Execute("70÷3=23, 1", $true, $false, $false, $false, $false, $true, 1, $false, "73÷8=9, 1", 2)
if (-not $found) { throw "Replacement failed: '70÷3=23, 1' -> '73÷8=9, 1'" }
$found = $d.Content.Find.Execute("31÷5=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "37÷2=18, 1", 2)
if (-not $found) { throw "Replacement failed: '31÷5=6, 1' -> '37÷2=18, 1'" }
$found = $d.Content.Find.Execute("38÷8=4, 6", $true, $false, $false, $false, $false, $true, 1, $false, "60÷8=7, 4", 2)
if (-not $found) { throw "Replacement failed: '38÷8=4, 6' -> '60÷8=7, 4'" }
$found = $d.Content.Find.Execute("62÷3=20, 2", $true, $false, $false, $false, $false, $true, 1, $false, "53÷2=26, 1", 2)
if (-not $found) { throw "Replacement failed: '62÷3=20, 2' -> '53÷2=26, 1'" }
$found = $d.Content.Find.Execute("76÷4=19, 0", $true, $false, $false, $false, $false, $true, 1, $false, "97÷8=12, 1", 2)
if (-not $found) { throw "Replacement failed: '76÷4=19, 0' -> '97÷8=12, 1'" }
$found = $d.Content.Find.Execute("91÷2=45, 1", $true, $false, $false, $false, $false, $true, 1, $false, "32÷2=16, 0", 2)
if (-not $found) { throw "Replacement failed: '91÷2=45, 1' -> '32÷2=16, 0'" }
$found = $d.Content.Find.Execute("24÷6=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "51÷3=17, 0", 2)
if (-not $found) { throw "Replacement failed: '24÷6=4, 0' -> '51÷3=17, 0'" }
$found = $d.Content.Find.Execute("92÷8=11, 4", $true, $false, $false, $false, $false, $true, 1, $false, "73÷7=10, 3", 2)
if (-not $found) { throw "Replacement failed: '92÷8=11, 4' -> '73÷7=10, 3'" }
$found = $d.Content.Find.Execute("84÷6=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "97÷5=19, 2", 2)
if (-not $found) { throw "Replacement failed: '84÷6=14, 0' -> '97÷5=19, 2'" }
